$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B29").Value = "H10"
$ws.Range("C37").Select()
